$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows before row 156; this shifts the former rows 156-236
# down to 159-239 (and carries column D's date-number format along).
$ws.Rows("156:158").Insert()

# Row 156 (new)
$ws.Range("A156").Value = 3
$ws.Range("B156").Value = "Femacal de La Calera"
$ws.Range("C156").Value = "Coquimbo"
$ws.Range("D156").Value = 44839
$ws.Range("E156").Value = 5
$ws.Range("F156").Value = "Fruta"
$ws.Range("G156").Value = 100101
$ws.Range("H156").Value = "Berries"
$ws.Range("I156").Value = 100101001
$ws.Range("J156").Value = "Arándano (blue)"
$ws.Range("K156").Value = "Sin especificar"
$ws.Range("L156").Value = "Primera"
$ws.Range("M156").Value = 45
$ws.Range("N156").Value = 10000
$ws.Range("O156").Value = 10000
$ws.Range("P156").Value = 10000
$ws.Range("Q156").Value = "$/bandeja 12 canastillos 125 gramos"
$ws.Range("R156").Value = "Provincia de Limarí"
$ws.Range("S156").Value = 6667
$ws.Range("T156").Value = 1.5

# Row 157 (new)
$ws.Range("A157").Value = 3
$ws.Range("B157").Value = "Femacal de La Calera"
$ws.Range("C157").Value = "Coquimbo"
$ws.Range("D157").Value = 44839
$ws.Range("E157").Value = 5
$ws.Range("F157").Value = "Fruta"
$ws.Range("G157").Value = 100101
$ws.Range("H157").Value = "Berries"
$ws.Range("I157").Value = 100101001
$ws.Range("J157").Value = "Arándano (blue)"
$ws.Range("K157").Value = "Sin especificar"
$ws.Range("L157").Value = "Primera"
$ws.Range("M157").Value = 25
$ws.Range("N157").Value = 12000
$ws.Range("O157").Value = 12000
$ws.Range("P157").Value = 12000
$ws.Range("Q157").Value = "$/bandeja 2 kilos"
$ws.Range("R157").Value = "Provincia de Quillota"
$ws.Range("S157").Value = 6000
$ws.Range("T157").Value = 2

# Row 158 (new)
$ws.Range("A158").Value = 3
$ws.Range("B158").Value = "Femacal de La Calera"
$ws.Range("C158").Value = "Coquimbo"
$ws.Range("D158").Value = 44839
$ws.Range("E158").Value = 5
$ws.Range("F158").Value = "Fruta"
$ws.Range("G158").Value = 100101
$ws.Range("H158").Value = "Berries"
$ws.Range("I158").Value = 100101001
$ws.Range("J158").Value = "Arándano (blue)"
$ws.Range("K158").Value = "Sin especificar"
$ws.Range("L158").Value = "Segunda"
$ws.Range("M158").Value = 20
$ws.Range("N158").Value = 10000
$ws.Range("O158").Value = 10000
$ws.Range("P158").Value = 10000
$ws.Range("Q158").Value = "$/bandeja 2 kilos"
$ws.Range("R158").Value = "Provincia de Quillota"
$ws.Range("S158").Value = 5000
$ws.Range("T158").Value = 2
